# Add "ZincCoat" component (with HHV-consistent impact data) as a new row
# (row 20) on every sheet of the impact_items workbook: info, GWP,
# H_Ecosystems, H_Health, H_Resources.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# info sheet: ID / functional_unit columns only
# ---------------------------------------------------------------------
$wsInfo = $wb.Worksheets.Item("info")
$wsInfo.Range("A20").Value = "ZincCoat"
$wsInfo.Range("B20").Value = "m2"

# ---------------------------------------------------------------------
# GWP sheet: fix a pre-existing reference typo on row 19 (ReCiPe ->
# TRACI, matching every other row), then append the ZincCoat row.
# ---------------------------------------------------------------------
$wsGWP = $wb.Worksheets.Item("GWP")
$wsGWP.Range("G19").Value = "ecoinvent 3.8 - cutoff, TRACI"

$wsGWP.Range("A20").Value = "ZincCoat"
$wsGWP.Range("B20").Value = "kg CO2-eq"
$wsGWP.Range("C20").Value = 5.4322602
$wsGWP.Range("D20").Formula = "=C20*0.9"
$wsGWP.Range("E20").Formula = "=C20*1.1"
$wsGWP.Range("F20").Value = "uniform"
$wsGWP.Range("G20").Value = "ecoinvent 3.8 - cutoff, TRACI"

# ---------------------------------------------------------------------
# H_Ecosystems sheet
# ---------------------------------------------------------------------
$wsEco = $wb.Worksheets.Item("H_Ecosystems")
$wsEco.Range("A20").Value = "ZincCoat"
$wsEco.Range("B20").Value = "points"
$wsEco.Range("C20").Value = 0.12593618000000001
$wsEco.Range("D20").Formula = "=C20*0.9"
$wsEco.Range("E20").Formula = "=C20*1.1"
$wsEco.Range("F20").Value = "uniform"
$wsEco.Range("G20").Value = "ecoinvent 3.8 - cutoff, ReCiPe Endpoint (H,A)"

# ---------------------------------------------------------------------
# H_Health sheet
# ---------------------------------------------------------------------
$wsHealth = $wb.Worksheets.Item("H_Health")
$wsHealth.Range("A20").Value = "ZincCoat"
$wsHealth.Range("B20").Value = "points"
$wsHealth.Range("C20").Value = 0.39828018999999998
$wsHealth.Range("D20").Formula = "=C20*0.9"
$wsHealth.Range("E20").Formula = "=C20*1.1"
$wsHealth.Range("F20").Value = "uniform"
$wsHealth.Range("G20").Value = "ecoinvent 3.8 - cutoff, ReCiPe Endpoint (H,A)"

# ---------------------------------------------------------------------
# H_Resources sheet
# ---------------------------------------------------------------------
$wsRes = $wb.Worksheets.Item("H_Resources")
$wsRes.Range("A20").Value = "ZincCoat"
$wsRes.Range("B20").Value = "points"
$wsRes.Range("C20").Value = 0.33182298999999998
$wsRes.Range("D20").Formula = "=C20*0.9"
$wsRes.Range("E20").Formula = "=C20*1.1"
$wsRes.Range("F20").Value = "uniform"
$wsRes.Range("G20").Value = "ecoinvent 3.8 - cutoff, ReCiPe Endpoint (H,A)"

# ---------------------------------------------------------------------
# Leave the cursor/selection on each sheet where the author left it, and
# finish with "info" as the active tab (matches the saved view state).
# ---------------------------------------------------------------------
[void]$wsGWP.Range("A20").Select()
[void]$wsEco.Range("D23").Select()
[void]$wsHealth.Range("D22").Select()
[void]$wsRes.Range("A20").Select()
[void]$wsInfo.Range("C22").Select()
[void]$wsInfo.Activate()
